$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.002") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values and mangles the exact displayed text (e.g. drops trailing zeros).
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D16",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D39",
    "D40",
    "D41",
    "D42",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.027.83"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.648.57"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "218.02"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "0.5210"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.2611"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").Value = "0.06266"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "20.44"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "4.462"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.613.18"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "0.5429"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "0.0₅8068"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "26.036.51"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "4.557"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").Value = "191.37"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "10.02"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "5.973"
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "138.55"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "0.1229"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "7.231"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "16.12"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "1.399"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "3.490"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "3.227"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").Value = "1.518"
$ws.Range("E33").Value = "  -8.29%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "0.9431"
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").Value = "2.744"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "0.5707"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "5.844"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "0.8438"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "100.47"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "1.001.40"
$ws.Range("E43").Value = "  -4.49%  "
$ws.Range("D44").Value = "1.792.54"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈107"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "56.52"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "0.4291"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "7.877"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "0.05149"
$ws.Range("D51").Value = "1.465"
$ws.Range("E51").Value = "  -1.42%  "

# Restore the original (default) style on the forced cells so only the
# value changes and no stray formatting is introduced.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
